# Sex_determination_outline.pptx - "everything uploaded to plos site,
# all figures passed pace system"
#
# The only substantive content change in this commit is on the sex
# determination diagram slide: the small Greek "rho" (ρ) label used
# for the recombination-rate annotation is no longer bold (it stays
# italic). Find that textbox and turn Bold off.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "rho" textbox. It lives a couple of levels deep inside a
# group on the slide, so walk the shape tree (including nested groups)
# looking for the run of text that is exactly the Greek letter rho.
$rho = [char]0x03C1
$target = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq $rho) {
        $target = $shp.TextFrame.TextRange
    }

    if ($shp.Type -eq 6) {
        # msoGroup - look one level into the group's members.
        $items = $shp.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            $sub = $items.Item($j)
            if ($sub.HasTextFrame -eq -1 -and $sub.TextFrame.TextRange.Text -eq $rho) {
                $target = $sub.TextFrame.TextRange
            }
        }
    }
}

if ($target -ne $null) {
    # Turn off Bold (italic / size / language stay untouched).
    $target.Font.Bold = $false
}
